# Update odds values for the week of 2025-05-08 (FlashScore export)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Range("G7").Value = 1.52
$ws.Range("H7").Value = 3.75
$ws.Range("J7").Value = 1.8
$ws.Range("K7").Value = 1.8
$ws.Range("L7").Value = 1.37
$ws.Range("M7").Value = 2.5
$ws.Range("N7").Value = 5.6
$ws.Range("O7").Value = 5.9
$ws.Range("P7").Value = 6.9
$ws.Range("Q7").Value = 8.75
$ws.Range("R7").Value = 10.25
$ws.Range("S7").Value = 22
$ws.Range("T7").Value = 10
$ws.Range("U7").Value = 6.5
$ws.Range("V7").Value = 14.5
$ws.Range("W7").Value = 65
$ws.Range("X7").Value = 450
$ws.Range("Y7").Value = 11.5
$ws.Range("Z7").Value = 26
$ws.Range("AA7").Value = 14.5
$ws.Range("AC7").Value = 45
$ws.Range("AD7").Value = 45

# Row 10
$ws.Range("G10").Value = 1.5
$ws.Range("H10").Value = 4.5
$ws.Range("I10").Value = 5.5
$ws.Range("Z10").Value = 29
$ws.Range("AB10").Value = 51
$ws.Range("AF10").Value = 17

# Row 12
$ws.Range("G12").Value = 2.27
$ws.Range("H12").Value = 3.05
$ws.Range("J12").Value = 2.32
$ws.Range("L12").Value = 1.52
$ws.Range("M12").Value = 2.22
$ws.Range("N12").Value = 6
$ws.Range("R12").Value = 22
$ws.Range("V12").Value = 18.5
$ws.Range("Y12").Value = 7.3
$ws.Range("AH12").Value = 2.35
$ws.Range("AI12").Value = 2.02
$ws.Range("AJ12").Value = 1.62

# Row 13
$ws.Range("G13").Value = 1.47
$ws.Range("H13").Value = 3.8
$ws.Range("I13").Value = 7.2
$ws.Range("J13").Value = 2
$ws.Range("K13").Value = 1.72
$ws.Range("O13").Value = 6
$ws.Range("P13").Value = 8.25
$ws.Range("Q13").Value = 9.5
$ws.Range("U13").Value = 7.5
$ws.Range("Y13").Value = 16
$ws.Range("Z13").Value = 45
$ws.Range("AA13").Value = 22
$ws.Range("AB13").Value = 175
$ws.Range("AC13").Value = 90
$ws.Range("AD13").Value = 90
$ws.Range("AG13").Value = 1.34
$ws.Range("AH13").Value = 3

# Row 19
$ws.Range("J19").Value = 2.1
$ws.Range("K19").Value = 1.7
$ws.Range("L19").Value = 1.44
$ws.Range("M19").Value = 2.63
$ws.Range("O19").Value = 5.5
$ws.Range("P19").Value = 9.5
$ws.Range("R19").Value = 15
$ws.Range("S19").Value = 41
$ws.Range("T19").Value = 7.5
$ws.Range("V19").Value = 26
$ws.Range("W19").Value = 101
$ws.Range("AA19").Value = 29
$ws.Range("AD19").Value = 81
$ws.Range("AG19").Value = 1.33
$ws.Range("AH19").Value = 3.25
$ws.Range("AI19").Value = 2.5
$ws.Range("AJ19").Value = 1.5

# Row 20
$ws.Range("N20").Value = 5
$ws.Range("U20").Value = 8
$ws.Range("W20").Value = 101
$ws.Range("AE20").Value = 1.08
$ws.Range("AF20").Value = 7.5
$ws.Range("AI20").Value = 2.5
$ws.Range("AJ20").Value = 1.5

# Row 21
$ws.Range("G21").Value = 1.45
$ws.Range("H21").Value = 3.4
$ws.Range("I21").Value = 7.5
$ws.Range("K21").Value = 1.58
$ws.Range("N21").Value = 5
$ws.Range("Q21").Value = 9.5
$ws.Range("T21").Value = 7
$ws.Range("W21").Value = 101
$ws.Range("Y21").Value = 15

# Row 22
$ws.Range("G22").Value = 6
$ws.Range("H22").Value = 3.6
$ws.Range("I22").Value = 1.53
$ws.Range("J22").Value = 1.95
$ws.Range("K22").Value = 1.85
$ws.Range("L22").Value = 1.36
$ws.Range("M22").Value = 3
$ws.Range("N22").Value = 15
$ws.Range("T22").Value = 9.5
$ws.Range("U22").Value = 7.5
$ws.Range("V22").Value = 17
$ws.Range("W22").Value = 51
$ws.Range("X22").Value = 351
$ws.Range("Y22").Value = 6.5
$ws.Range("Z22").Value = 7
$ws.Range("AA22").Value = 8.5
$ws.Range("AC22").Value = 13
$ws.Range("AD22").Value = 29
$ws.Range("AE22").Value = 1.04
$ws.Range("AF22").Value = 9.5
$ws.Range("AG22").Value = 1.25
$ws.Range("AH22").Value = 3.5
$ws.Range("AI22").Value = 1.91
$ws.Range("AJ22").Value = 1.8

# Row 23
$ws.Range("K23").Value = 1.41
$ws.Range("R23").Value = 23
$ws.Range("AA23").Value = 13
$ws.Range("AC23").Value = 34
$ws.Range("AE23").Value = 1.08
$ws.Range("AG23").Value = 1.5
$ws.Range("AH23").Value = 2.37

# Row 24
$ws.Range("G24").Value = 1.58

# Row 25
$ws.Range("G25").Value = 2.05
$ws.Range("I25").Value = 3.4
$ws.Range("J25").Value = 2.7
$ws.Range("K25").Value = 1.44
$ws.Range("N25").Value = 5.5
$ws.Range("P25").Value = 10
$ws.Range("S25").Value = 41
$ws.Range("T25").Value = 6.5
$ws.Range("U25").Value = 6.5
$ws.Range("V25").Value = 21
$ws.Range("W25").Value = 81
$ws.Range("Y25").Value = 7.5
$ws.Range("AA25").Value = 15
$ws.Range("AC25").Value = 41
$ws.Range("AD25").Value = 51
